{"js": "// Apply the tracked-changes described in the commit:\n// \"Chinh sua PHPWord, th\u00eam \u0111\u00e1nh d\u1ea5u tag script\"\n//\n// 1. Update the opening date/time (\"H\u1ed3i 14 gi\u1edd 43 ph\u00fat, ng\u00e0y 16 ...\")\n// 2. Update the \"Quy\u1ebft \u0111\u1ecbnh tr\u01b0ng c\u1ea7u gi\u00e1m \u0111\u1ecbnh s\u1ed1 ...\" reference\n// 3. Update the issuing authority name (2 occurrences)\n// 4. Update the \"B\u00ean giao\" signer name (2 occurrences)\n// 5. Replace the three phone/IMEI paragraphs with a single video paragraph\n// 6. Update the closing date/time (\"Vi\u1ec7c giao, nh\u1eadn k\u1ebft th\u00fac h\u1ed3i 14 gi\u1edd 58 ph\u00fat, ng\u00e0y 16 ...\")\n\nconst body = context.document.body;\n\n// ---- 1) Opening timestamp: \"H\u1ed3i 14 gi\u1edd 43 ph\u00fat, ng\u00e0y 16 th\u00e1ng 05 n\u0103m 2023\" ----\nconst hoiResults = body.search(\"H\u1ed3i\", { matchCase: true });\nhoiResults.load(\"text\");\nawait context.sync();\n\n// Use the paragraph containing \"H\u1ed3i\" (start of the bi\u00ean b\u1ea3n) to scope the\n// replacements so we only touch the opening sentence, not the closing one.\nconst openParaRange = hoiResults.items[0].paragraphs.getFirst();\n\nconst openHour = openParaRange.search(\" 14\", { matchCase: true });\nconst openMinute = openParaRange.search(\"43\", { matchCase: true });\nconst openDay = openParaRange.search(\"16\", { matchCase: true });\nopenHour.load(\"text\");\nopenMinute.load(\"text\");\nopenDay.load(\"text\");\nawait context.sync();\n\nopenHour.items[0].insertText(\" 15\", \"Replace\");\nopenMinute.items[0].insertText(\"08\", \"Replace\");\nopenDay.items[0].insertText(\"19\", \"Replace\");\nawait context.sync();\n\n// ---- 2) Decision number & date ----\nconst decisionNoResults = body.search(\"59/Q\u0110-CS\u0110T\", { matchCase: true });\ndecisionNoResults.load(\"text\");\nawait context.sync();\ndecisionNoResults.items[0].insertText(\"185/Q\u0110-CQ\u0110T\", \"Replace\");\nawait context.sync();\n\nconst decisionDateResults = body.search(\"15/05/2023 \", { matchCase: true });\ndecisionDateResults.load(\"text\");\nawait context.sync();\ndecisionDateResults.items[0].insertText(\"18/05/2023 \", \"Replace\");\nawait context.sync();\n\n// ---- 3) Issuing authority name (appears twice, same replacement both times) ----\nconst authorityResults = body.search(\n  \"C\u01a1 quan C\u1ea3nh s\u00e1t \u0111i\u1ec1u tra C\u00f4ng an huy\u1ec7n Qu\u1ea3ng H\u00f2a\",\n  { matchCase: true }\n);\nauthorityResults.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < authorityResults.items.length; i++) {\n  authorityResults.items[i].insertText(\n    \"C\u01a1 quan C\u1ea3nh s\u00e1t \u0111i\u1ec1u tra C\u00f4ng an th\u00e0nh ph\u1ed1 Cao B\u1eb1ng\",\n    \"Replace\"\n  );\n}\nawait context.sync();\n\n// ---- 4) \"B\u00ean giao\" signer name (inline text + signature table, same replacement) ----\nconst signerResults = body.search(\"Vi M\u1ea1nh Hi\u1ec3n\", { matchCase: true });\nsignerResults.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < signerResults.items.length; i++) {\n  signerResults.items[i].insertText(\"H\u00e0 Th\u1ebf Duy\", \"Replace\");\n}\nawait context.sync();\n\n// ---- 5) Replace the 3 phone/IMEI paragraphs with a single video paragraph ----\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst phoneParaIndexes = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"\u0111i\u1ec7n tho\u1ea1i di \u0111\u1ed9ng\") !== -1) {\n    phoneParaIndexes.push(i);\n  }\n}\n\nconst newExhibitText =\n  '- 01 \u0111o\u1ea1n video c\u00f3 t\u00ean \\u201Cch0_20230416101915_013.mp4\\u201D, ' +\n  'm\u00e3 MD5: 8bc02026d02c347db852c64d425468d2 \u0111\u01b0\u1ee3c l\u01b0u tr\u1eef trong th\u1ebb nh\u1edb ' +\n  'c\u00f3 ch\u1eef Pioneer 32GB, ni\u00eam phong trong 01 b\u00ec th\u01b0 ghi \"Th\u1ebb nh\u1edb camera ' +\n  'h\u00e0nh tr\u00ecnh xe \u00f4 t\u00f4 BKS: 000.46';\n\nparagraphs.items[phoneParaIndexes[0]].insertText(newExhibitText, \"Replace\");\nparagraphs.items[phoneParaIndexes[1]].delete();\nparagraphs.items[phoneParaIndexes[2]].delete();\nawait context.sync();\n\n// ---- 6) Closing timestamp: \"k\u1ebft th\u00fac h\u1ed3i 14 gi\u1edd 58 ph\u00fat, ng\u00e0y 16 th\u00e1ng 05 n\u0103m 2023\" ----\nconst endResults = body.search(\"k\u1ebft th\u00fac h\u1ed3i\", { matchCase: true });\nendResults.load(\"text\");\nawait context.sync();\nconst endParaRange = endResults.items[0].paragraphs.getFirst();\n\nconst endHour = endParaRange.search(\"14\", { matchCase: true });\nconst endMinute = endParaRange.search(\"58\", { matchCase: true });\nconst endDay = endParaRange.search(\"16 \", { matchCase: true });\nendHour.load(\"text\");\nendMinute.load(\"text\");\nendDay.load(\"text\");\nawait context.sync();\n\nendHour.items[0].insertText(\"15\", \"Replace\");\nendMinute.items[0].insertText(\"23\", \"Replace\");\nendDay.items[0].insertText(\"19 \", \"Replace\");\nawait context.sync();\n", "ps1": "# Apply the tracked-changes described in the commit:\n# \"Chinh sua PHPWord, th\u00eam \u0111\u00e1nh d\u1ea5u tag script\"\n#\n# 1. Update the opening date/time (\"H\u1ed3i 14 gi\u1edd 43 ph\u00fat, ng\u00e0y 16 ...\")\n# 2. Update the \"Quy\u1ebft \u0111\u1ecbnh tr\u01b0ng c\u1ea7u gi\u00e1m \u0111\u1ecbnh s\u1ed1 ...\" reference\n# 3. Update the issuing authority name (2 occurrences)\n# 4. Update the \"B\u00ean giao\" signer name (2 occurrences)\n# 5. Replace the three phone/IMEI paragraphs with a single video paragraph\n# 6. Update the closing date/time (\"Vi\u1ec7c giao, nh\u1eadn k\u1ebft th\u00fac h\u1ed3i 14 gi\u1edd 58 ph\u00fat, ng\u00e0y 16 ...\")\n\n$d = $word.ActiveDocument\n\nfunction Replace-AllText($searchText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $searchText\n    $find.Replacement.Text = $replaceText\n    $find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n}\n\nfunction Replace-TextInRange($range, $searchText, $replaceText) {\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $searchText\n    $find.Replacement.Text = $replaceText\n    $find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n}\n\n# ---- Locate the opening & closing timestamp paragraphs up-front (their\n# ---- index would shift once we delete paragraphs further down, so resolve\n# ---- the time edits before touching the exhibit list). ----\n$paraCount = $d.Paragraphs.Count\n$openParaIdx = -1\n$closeParaIdx = -1\nfor ($i = 1; $i -le $paraCount; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($openParaIdx -eq -1 -and $t.StartsWith(\"H\u1ed3i\")) {\n        $openParaIdx = $i\n    }\n    if ($t.Contains(\"k\u1ebft th\u00fac h\u1ed3i\")) {\n        $closeParaIdx = $i\n    }\n}\n\n# ---- 1) Opening timestamp: \"H\u1ed3i 14 gi\u1edd 43 ph\u00fat, ng\u00e0y 16 th\u00e1ng 05 n\u0103m 2023\" ----\n$openRange = $d.Paragraphs.Item($openParaIdx).Range\nReplace-TextInRange $openRange \" 14\" \" 15\"\n$openRange = $d.Paragraphs.Item($openParaIdx).Range\nReplace-TextInRange $openRange \"43\" \"08\"\n$openRange = $d.Paragraphs.Item($openParaIdx).Range\nReplace-TextInRange $openRange \"16\" \"19\"\n\n# ---- 2) Decision number & date ----\nReplace-AllText \"59/Q\u0110-CS\u0110T\" \"185/Q\u0110-CQ\u0110T\"\nReplace-AllText \"15/05/2023 \" \"18/05/2023 \"\n\n# ---- 3) Issuing authority name (appears twice, same replacement both times) ----\nReplace-AllText \"C\u01a1 quan C\u1ea3nh s\u00e1t \u0111i\u1ec1u tra C\u00f4ng an huy\u1ec7n Qu\u1ea3ng H\u00f2a\" \"C\u01a1 quan C\u1ea3nh s\u00e1t \u0111i\u1ec1u tra C\u00f4ng an th\u00e0nh ph\u1ed1 Cao B\u1eb1ng\"\n\n# ---- 4) \"B\u00ean giao\" signer name (inline text + signature table, same replacement) ----\nReplace-AllText \"Vi M\u1ea1nh Hi\u1ec3n\" \"H\u00e0 Th\u1ebf Duy\"\n\n# ---- 5) Replace the 3 phone/IMEI paragraphs with a single video paragraph ----\n$paraCount = $d.Paragraphs.Count\n$phoneParaIdx = @()\nfor ($i = 1; $i -le $paraCount; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t.Contains(\"\u0111i\u1ec7n tho\u1ea1i di \u0111\u1ed9ng\")) {\n        $phoneParaIdx += $i\n    }\n}\n\n$newExhibitText = \"- 01 \u0111o\u1ea1n video c\u00f3 t\u00ean \" + [char]0x201C + \"ch0_20230416101915_013.mp4\" + [char]0x201D + \", m\u00e3 MD5: 8bc02026d02c347db852c64d425468d2 \u0111\u01b0\u1ee3c l\u01b0u tr\u1eef trong th\u1ebb nh\u1edb c\u00f3 ch\u1eef Pioneer 32GB, ni\u00eam phong trong 01 b\u00ec th\u01b0 ghi \" + [char]0x22 + \"Th\u1ebb nh\u1edb camera h\u00e0nh tr\u00ecnh xe \u00f4 t\u00f4 BKS: 000.46\"\n\n$d.Paragraphs.Item($phoneParaIdx[0]).Range.Text = $newExhibitText\n$d.Paragraphs.Item($phoneParaIdx[1]).Range.Delete() | Out-Null\n$d.Paragraphs.Item($phoneParaIdx[1]).Range.Delete() | Out-Null\n\n# ---- 6) Closing timestamp: \"k\u1ebft th\u00fac h\u1ed3i 14 gi\u1edd 58 ph\u00fat, ng\u00e0y 16 th\u00e1ng 05 n\u0103m 2023\" ----\n# Paragraph indices above $phoneParaIdx shifted down by 2 once the two extra\n# exhibit paragraphs were deleted, so re-locate the closing paragraph by text.\n$paraCount = $d.Paragraphs.Count\n$closeParaIdx = -1\nfor ($i = 1; $i -le $paraCount; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text.Contains(\"k\u1ebft th\u00fac h\u1ed3i\")) {\n        $closeParaIdx = $i\n        break\n    }\n}\n\n$closeRange = $d.Paragraphs.Item($closeParaIdx).Range\nReplace-TextInRange $closeRange \"14\" \"15\"\n$closeRange = $d.Paragraphs.Item($closeParaIdx).Range\nReplace-TextInRange $closeRange \"58\" \"23\"\n$closeRange = $d.Paragraphs.Item($closeParaIdx).Range\nReplace-TextInRange $closeRange \"16 \" \"19 \"\n"}
